# Rebuild the sheet content/row-heights to match the post-edit layout described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (old row 24) is removed entirely, shifting nothing else (it was the final row).
$ws.Rows(24).Delete()

# Row 1
$ws.Range("B1").Value = "'Ementa atual:"
$ws.Range("C1").Value = "'Ementa modificada (dados modificados em vermelho):"

# Row 2
$ws.Range("B2").Value = "'LOT2042"
$ws.Range("C2").Value = "'LOT2042"

# Row 3
$ws.Range("A3").Value = "'Nome:"
$ws.Range("B3").Value = "' Processos Bioquímicos"
$ws.Range("C3").Value = "' Processos Bioquímicos"

# Row 4
$ws.Range("A4").Value = "'Name:"
$ws.Range("B4").Value = "'Biochemical Processes"
$ws.Range("C4").Value = "'Biochemical Processes"

# Row 5
$ws.Range("A5").Value = "'Créditos-aula:"
$ws.Range("B5").Value = "'4"
$ws.Range("C5").Value = "'4"

# Row 6
$ws.Range("A6").Value = "'Créditos-trabalho"
$ws.Range("B6").Value = "'0"
$ws.Range("C6").Value = "'0"

# Row 7
$ws.Range("A7").Value = "'Carga horária:"
$ws.Range("B7").Value = "'60 h"
$ws.Range("C7").Value = "'60 h"

# Row 8
$ws.Range("A8").Value = "'Ativação:"
$ws.Range("B8").Value = "'01/01/2018"
$ws.Range("C8").Value = "'01/01/2018"

# Row 9
$ws.Range("A9").Value = "'Semestre ideal:"
$ws.Range("B9").Value = "'EQD-9,EQN-10"
$ws.Range("C9").Value = "'EQD-9,EQN-10"

# Row 10
$ws.Range("A10").Value = "'Objetivos:"
$ws.Range("B10").Value = "'1814052 - Silvio Silverio da Silva"
$ws.Range("C10").Value = "'1814052 - Silvio Silverio da Silva"

# Row 11
$ws.Range("A11").Value = "'Objectives:"

# Row 12
$ws.Range("A12").Value = "'Docentes responsáveis:"

# Row 13
$ws.Range("A13").Value = "'Programa resumido:"
$ws.Range("B13").Value = "'Semestral"
$ws.Range("C13").Value = "'Semestral"
$ws.Rows(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "'Short syllabus:"
$ws.Range("B14").Value = "'Introduction. Fermentative and enzymatic processes. Biochemical processes applied to food industry. Biochemical processes of industrial Importance. Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes. Variation of scale in bioprocesses. Introduction to techniques of separation/purification of biotechnological products."
$ws.Range("C14").Value = "'Introduction. Fermentative and enzymatic processes. Biochemical processes applied to food industry. Biochemical processes of industrial Importance. Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes. Variation of scale in bioprocesses. Introduction to techniques of separation/purification of biotechnological products."

# Row 15
$ws.Range("A15").Value = "'Programa:"
$ws.Range("B15").Value = "'01/01/2018"
$ws.Range("C15").Value = "'01/01/2018"
$ws.Rows(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "'Syllabus:"
$ws.Range("B16").Value = "'Introduction: importance of bioprocesses and industrial applications.Fermentative and enzymatic processes: types of fermentative processes, raw materials, products obtainment.Biochemical processes applied to food industry: food processing, phases of food products processing, biochemical alterations in food, lipids oxidation, enzymatic and not enzymatic darkness, industrial controls of biochemical alterations.Biotechnological processes of industrial importance: description and study of cases of some biotechnological processes.Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes: transfer of oxygen and microbial respiration: transfer of mass (convection transfer in gas-liquid system, microbial respiration, transfer of oxygen from the gas bubble to the cell); Transfer of O2 in the bioreactor (Effects of sizing and operating aspects of the bioreactor - bubbles, aeration, agitation and medium properties, antifoaming agents, temperature, gas pressure and oxygen partial pressure). Transfer of power and oxygen in agitated and aerated bioreactor. Scale variation in bioprocesses.Introduction to separation/purification of biotechnological products"
$ws.Range("C16").Value = "'Introduction: importance of bioprocesses and industrial applications.Fermentative and enzymatic processes: types of fermentative processes, raw materials, products obtainment.Biochemical processes applied to food industry: food processing, phases of food products processing, biochemical alterations in food, lipids oxidation, enzymatic and not enzymatic darkness, industrial controls of biochemical alterations.Biotechnological processes of industrial importance: description and study of cases of some biotechnological processes.Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes: transfer of oxygen and microbial respiration: transfer of mass (convection transfer in gas-liquid system, microbial respiration, transfer of oxygen from the gas bubble to the cell); Transfer of O2 in the bioreactor (Effects of sizing and operating aspects of the bioreactor - bubbles, aeration, agitation and medium properties, antifoaming agents, temperature, gas pressure and oxygen partial pressure). Transfer of power and oxygen in agitated and aerated bioreactor. Scale variation in bioprocesses.Introduction to separation/purification of biotechnological products"

# Row 17
$ws.Range("A17").Value = "'Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows(17).AutoFit()

# Row 18
$ws.Range("A18").Value = "'Método:"
$ws.Range("B18").Value = "'1814052 - Silvio Silverio da Silva"
$ws.Range("C18").Value = "'1814052 - Silvio Silverio da Silva"
$ws.Rows(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "'Critério:"
$ws.Range("B19").Value = "'Os alunos serão avaliados formalmente por uma prova teórica (P) e trabalhos (T). A ponderação das notas será de 70% para a prova teórica (P) e 30% para a média aritmética das notas dos trabalhos (T), ou seja: Média Final do período letivo normal (MF) = (0,7xP +0,3xT)."
$ws.Range("C19").Value = "'Os alunos serão avaliados formalmente por uma prova teórica (P) e trabalhos (T). A ponderação das notas será de 70% para a prova teórica (P) e 30% para a média aritmética das notas dos trabalhos (T), ou seja: Média Final do período letivo normal (MF) = (0,7xP +0,3xT)."

# Row 20
$ws.Range("A20").Value = "'Norma de recuperação:"
$ws.Range("B20").Value = "'Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Range("C20").Value = "'Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."

# Row 21
$ws.Range("A21").Value = "'Bibliografia:"
$ws.Range("B21").Value = "'A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0."
$ws.Range("C21").Value = "'A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0."
$ws.Rows(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "'Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows(22).AutoFit()

# Row 23
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "'LOT2041 -  Engenharia Bioquímica  (Requisito fraco)`n"
$ws.Range("C23").Value = "'LOT2041 -  Engenharia Bioquímica  (Requisito fraco)`n"
$ws.Rows(23).RowHeight = 30
